$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect before editing, re-protect afterward.
$ws.Unprotect()

# --- Update the confidential disclaimer date in A44 (shared string) ---
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."
$ws.Cells.Item(44, 1).Value = $disclaimer

# --- Update Weight (col D) and Percent Change (col E) values for rows 2-41 ---
$ws.Cells.Item(2, 4).Value = 0.07795619663159882
$ws.Cells.Item(2, 5).Value = 0.002462271644162017
$ws.Cells.Item(3, 4).Value = 0.06976709557956315
$ws.Cells.Item(3, 5).Value = -0.004858072027943861
$ws.Cells.Item(4, 4).Value = 0.05421227060433081
$ws.Cells.Item(4, 5).Value = 0.001522881291403211
$ws.Cells.Item(5, 4).Value = 0.04757055833372648
$ws.Cells.Item(5, 5).Value = -0.0009018418026918162
$ws.Cells.Item(6, 4).Value = 0.04303055856960891
$ws.Cells.Item(6, 5).Value = -0.007030334591850207
$ws.Cells.Item(7, 4).Value = 0.03875076661791763
$ws.Cells.Item(7, 5).Value = -0.001780496712929125
$ws.Cells.Item(8, 4).Value = 0.03855026654715291
$ws.Cells.Item(8, 5).Value = -0.0002447531053051399
$ws.Cells.Item(9, 4).Value = 0.03494465608340803
$ws.Cells.Item(9, 5).Value = 0.004805278634440135
$ws.Cells.Item(10, 4).Value = 0.03271335566353729
$ws.Cells.Item(10, 5).Value = -0.004371417240508935
$ws.Cells.Item(11, 4).Value = 0.02778989951408218
$ws.Cells.Item(11, 5).Value = 0.006525198938992061
$ws.Cells.Item(12, 4).Value = 0.03168490824173232
$ws.Cells.Item(12, 5).Value = -0.002763819095477404
$ws.Cells.Item(13, 4).Value = 0.03187906897202435
$ws.Cells.Item(13, 5).Value = -0.008578550386844408
$ws.Cells.Item(14, 4).Value = 0.0271500684059065
$ws.Cells.Item(14, 5).Value = 0.00339378801042578
$ws.Cells.Item(15, 4).Value = 0.03017452351747889
$ws.Cells.Item(15, 5).Value = 0.0003664345914256018
$ws.Cells.Item(16, 4).Value = 0.02728717507194415
$ws.Cells.Item(16, 5).Value = -0.01569506726457404
$ws.Cells.Item(17, 4).Value = 0.02789899514082182
$ws.Cells.Item(17, 5).Value = -0.005136334812935961
$ws.Cells.Item(18, 4).Value = 0.02354342595650328
$ws.Cells.Item(18, 5).Value = -0.01472798316801927
$ws.Cells.Item(19, 4).Value = 0.02037581969146577
$ws.Cells.Item(19, 5).Value = 0.0008248317777295178
$ws.Cells.Item(20, 4).Value = 0.02181322828702175
$ws.Cells.Item(20, 5).Value = 0.0007299270072993469
$ws.Cells.Item(21, 4).Value = 0.020823996320234
$ws.Cells.Item(21, 5).Value = -0.002831858407079668
$ws.Cells.Item(22, 4).Value = 0.02169749846676417
$ws.Cells.Item(22, 5).Value = 0.004076779344317938
$ws.Cells.Item(23, 4).Value = 0.02024254611501628
$ws.Cells.Item(23, 5).Value = 0.007195606892633988
$ws.Cells.Item(24, 4).Value = 0.01992086144265698
$ws.Cells.Item(24, 5).Value = 0.008791924454575106
$ws.Cells.Item(25, 4).Value = 0.01767054300136812
$ws.Cells.Item(25, 5).Value = 0.0009177373602538008
$ws.Cells.Item(26, 4).Value = 0.01792824456290985
$ws.Cells.Item(26, 5).Value = 0.01509769094138558
$ws.Cells.Item(27, 4).Value = 0.01931656012643299
$ws.Cells.Item(27, 5).Value = -0.007345926349933252
$ws.Cells.Item(28, 4).Value = 0.01694830046704723
$ws.Cells.Item(28, 5).Value = -0.008898669983733787
$ws.Cells.Item(29, 4).Value = 0.01816294758692268
$ws.Cells.Item(29, 5).Value = -0.01571428571428568
$ws.Cells.Item(30, 4).Value = 0.01765049299429165
$ws.Cells.Item(30, 5).Value = -0.006339581036383768
$ws.Cells.Item(31, 4).Value = 0.01858635655989055
$ws.Cells.Item(31, 5).Value = 0.001618122977346204
$ws.Cells.Item(32, 4).Value = 0.01564785111100628
$ws.Cells.Item(32, 5).Value = -0.007160354249105105
$ws.Cells.Item(33, 4).Value = 0.01688446478275228
$ws.Cells.Item(33, 5).Value = 0.003562447611064501
$ws.Cells.Item(34, 4).Value = 0.008248513940651979
$ws.Cells.Item(34, 5).Value = -0.009008042895442325
$ws.Cells.Item(35, 4).Value = 0.0079709097985564
$ws.Cells.Item(35, 5).Value = 0.007139290140011578
$ws.Cells.Item(36, 4).Value = 0.007401401141670992
$ws.Cells.Item(36, 5).Value = 0.008644729503625159
$ws.Cells.Item(37, 4).Value = 0.006495317733641553
$ws.Cells.Item(37, 5).Value = 0.003631576558173366
$ws.Cells.Item(38, 4).Value = 0.007201490776996746
$ws.Cells.Item(38, 5).Value = -0.009437438584998348
$ws.Cells.Item(39, 4).Value = 0.007247782705099779
$ws.Cells.Item(39, 5).Value = -0.0005695455839876962
$ws.Cells.Item(40, 4).Value = 0.006861082936264566
$ws.Cells.Item(40, 5).Value = 0.007198263821740936
$ws.Cells.Item(41, 5).Value = -0.001452667240647298

# Restore sheet protection
$ws.Protect()
